$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sits right before the
#    document title ("Playing Card Run Specification").
# ------------------------------------------------------------------
$titleBm = $d.Bookmarks("_GoBack")
$titleBm.Delete()

# ------------------------------------------------------------------
# 2. Rewrite the "If you use clone the project from my github
#    address ..." sentence into the new, longer sentence.
# ------------------------------------------------------------------
$oldSentence = " you use clone the project from my github address to your local, please run under command to play game. Also you can change the number of players by change the parameter "
$newSentence = " you use fork my GitHub project address and clone it to your local, please run under command to play game in the forked project. Also you can change the number of players by change the parameter "

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $rng.Start
$rng.Text = $newSentence

# ------------------------------------------------------------------
# 3. Re-create the run boundaries the real edit produced.  Word
#    merges adjacent same-formatted runs together, but a bookmark
#    dropped at a given offset forces a split there that survives
#    even after the bookmark itself is removed again - use that to
#    reproduce every run boundary from the original edit.
# ------------------------------------------------------------------
$splitOffsets = @(8, 9, 13, 14, 17, 18, 20, 21, 24, 31, 32, 39, 52, 93, 105, 127)
foreach ($off in $splitOffsets) {
    $pos = $base + $off
    $splitRange = $d.Range($pos, $pos)
    $tmpName = "TempRunSplit"
    $tmpBm = $d.Bookmarks.Add($tmpName, $splitRange)
    $d.Bookmarks($tmpName).Delete()
}

# ------------------------------------------------------------------
# 4. Re-insert the "_GoBack" bookmark in its new home: splitting the
#    word "under" into "und" | bookmark | "er" right after "please
#    run und".
# ------------------------------------------------------------------
$bookmarkOffset = 82
$bmPos = $base + $bookmarkOffset
$insertionPoint = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $insertionPoint)
